$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Rename the header in C1 from "BSC/RNC" to "Technology"
$ws.Range("C1").Value = "Technology"

# Widen column C to fit the new header text
$ws.Columns.Item(3).ColumnWidth = 17

# Move the active selection to C1 (was B8)
$ws.Range("C1").Select()
